$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D33").Value = "C2905434"
$ws.Range("C28").Value = "PJ-3200"
$ws.Range("D28").Value = "C2689690"
$ws.Range("A28").Value = "Audio Jack"

$ws.Range("A28").Select()
